# This script reproduces the edit described by the diff: a new data row
# (for a later weekly Papaya price observation) was inserted as row 37 of
# the sheet, pushing the previously-existing rows 37-76 down to 38-77.
#
# The new row 37 contains:
#   Fecha (D) = 2022-04-18  (Excel serial 44669)
#   Volumen (M) = 50
#   Precio minimo/maximo/promedio (N/O/P) = 25000
#   Unidad de comercializacion (Q) = "$/bandeja 10 kilos"
#   Precio $/Kg (S) = 2500
#   Kg/unidad (T) = 10
# All other columns repeat the same constant values used by every row in
# this subset (Mercado, Region, Codreg, Tipo, Producto, Categoria, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 37; everything below (old rows 37-76)
# shifts down to rows 38-77, matching the dimension change A1:T76 -> A1:T77.
$ws.Rows.Item(37).Insert()

$ws.Range("A37").Value2 = 10
$ws.Range("B37").Value2 = "Vega Modelo de Temuco"
$ws.Range("C37").Value2 = "La Araucanía"
$ws.Range("D37").Value2 = 44669
$ws.Range("E37").Value2 = 9
$ws.Range("F37").Value2 = "Fruta"
$ws.Range("G37").Value2 = 100108
$ws.Range("H37").Value2 = "Tropicales y subtropicales"
$ws.Range("I37").Value2 = 100108004
$ws.Range("J37").Value2 = "Papaya"
$ws.Range("K37").Value2 = "Cultivar IV Región"
$ws.Range("L37").Value2 = "Primera"
$ws.Range("M37").Value2 = 50
$ws.Range("N37").Value2 = 25000
$ws.Range("O37").Value2 = 25000
$ws.Range("P37").Value2 = 25000
$ws.Range("Q37").Value2 = '$/bandeja 10 kilos'
$ws.Range("R37").Value2 = "Provincia del Elquí"
$ws.Range("S37").Value2 = 2500
$ws.Range("T37").Value2 = 10

# Make sure the new date cell keeps the same date/time number format as the
# rest of the "Fecha" column (style index 2 in the original workbook).
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
